$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows added below the existing table
$ws.Range("A10").Value = "magnitude and phase response"

$ws.Range("B13").Value = "StrainDifferential(mV)"
$ws.Range("B14").Value = "100Hz sin wave (-10mV,10mV)"

# New column width for column B (closest reachable quantized value to 17.62890625)
$ws.Columns.Item(2).ColumnWidth = 16.8

# Selection matches final diff state (B14 selected)
$ws.Range("B14").Select()
